# Extend the "T / Seno(T)" table from row 30 down to row 100.
# Column A keeps the same +0.01 step pattern that was already used for rows 2-30;
# column B keeps filling the SIN(A) formula down the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: literal step values for the 70 new rows (31-100) ---
$ws.Range("A31").Value = 0.3
$ws.Range("A32").Value = 0.31
$ws.Range("A33").Value = 0.32
$ws.Range("A34").Value = 0.33
$ws.Range("A35").Value = 0.34
$ws.Range("A36").Value = 0.35
$ws.Range("A37").Value = 0.36
$ws.Range("A38").Value = 0.37
$ws.Range("A39").Value = 0.38
$ws.Range("A40").Value = 0.38999999999999901
$ws.Range("A41").Value = 0.39999999999999902
$ws.Range("A42").Value = 0.40999999999999898
$ws.Range("A43").Value = 0.41999999999999899
$ws.Range("A44").Value = 0.42999999999999899
$ws.Range("A45").Value = 0.439999999999999
$ws.Range("A46").Value = 0.44999999999999901
$ws.Range("A47").Value = 0.45999999999999902
$ws.Range("A48").Value = 0.46999999999999897
$ws.Range("A49").Value = 0.47999999999999898
$ws.Range("A50").Value = 0.48999999999999899
$ws.Range("A51").Value = 0.499999999999999
$ws.Range("A52").Value = 0.50999999999999901
$ws.Range("A53").Value = 0.51999999999999902
$ws.Range("A54").Value = 0.52999999999999903
$ws.Range("A55").Value = 0.53999999999999904
$ws.Range("A56").Value = 0.54999999999999905
$ws.Range("A57").Value = 0.55999999999999905
$ws.Range("A58").Value = 0.56999999999999895
$ws.Range("A59").Value = 0.57999999999999896
$ws.Range("A60").Value = 0.58999999999999897
$ws.Range("A61").Value = 0.59999999999999898
$ws.Range("A62").Value = 0.60999999999999799
$ws.Range("A63").Value = 0.619999999999998
$ws.Range("A64").Value = 0.62999999999999801
$ws.Range("A65").Value = 0.63999999999999801
$ws.Range("A66").Value = 0.64999999999999802
$ws.Range("A67").Value = 0.65999999999999803
$ws.Range("A68").Value = 0.66999999999999804
$ws.Range("A69").Value = 0.67999999999999805
$ws.Range("A70").Value = 0.68999999999999795
$ws.Range("A71").Value = 0.69999999999999796
$ws.Range("A72").Value = 0.70999999999999797
$ws.Range("A73").Value = 0.71999999999999797
$ws.Range("A74").Value = 0.72999999999999798
$ws.Range("A75").Value = 0.73999999999999799
$ws.Range("A76").Value = 0.749999999999998
$ws.Range("A77").Value = 0.75999999999999801
$ws.Range("A78").Value = 0.76999999999999802
$ws.Range("A79").Value = 0.77999999999999803
$ws.Range("A80").Value = 0.78999999999999804
$ws.Range("A81").Value = 0.79999999999999805
$ws.Range("A82").Value = 0.80999999999999805
$ws.Range("A83").Value = 0.81999999999999695
$ws.Range("A84").Value = 0.82999999999999696
$ws.Range("A85").Value = 0.83999999999999697
$ws.Range("A86").Value = 0.84999999999999698
$ws.Range("A87").Value = 0.85999999999999699
$ws.Range("A88").Value = 0.869999999999997
$ws.Range("A89").Value = 0.87999999999999701
$ws.Range("A90").Value = 0.88999999999999702
$ws.Range("A91").Value = 0.89999999999999702
$ws.Range("A92").Value = 0.90999999999999703
$ws.Range("A93").Value = 0.91999999999999704
$ws.Range("A94").Value = 0.92999999999999705
$ws.Range("A95").Value = 0.93999999999999695
$ws.Range("A96").Value = 0.94999999999999696
$ws.Range("A97").Value = 0.95999999999999697
$ws.Range("A98").Value = 0.96999999999999698
$ws.Range("A99").Value = 0.97999999999999698
$ws.Range("A100").Value = 0.98999999999999699

# --- Column B: fill the =SIN(A#) formula down through the new rows ---
# First batch continues/extends the existing shared formula run (rows 31-66),
# second batch covers the remaining new rows (67-100).
$ws.Range("B31:B66").Formula = "=SIN(A31)"
$ws.Range("B67:B100").Formula = "=SIN(A67)"

# --- Update the view so it matches where the user ended up scrolled to/selecting ---
$ws.Range("D88").Select()
